$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "done changes in payment" - row 2 holds the payment-test credentials/url.
# B2 used to hold a stray numeric placeholder (123); it should hold the
# account password instead, and C2's login URL should point at the
# baramatimc test endpoint (matching row 3's site) rather than panvelmc.
$ws.Range("B2").Value = "User@12345"
$ws.Range("C2").Value = "http://testbaramatimc.ptaxcollection.com:8080/Pages/Login.aspx"

# Move / record the active selection at C2 (matches the saved view state).
$ws.Range("C2").Select() | Out-Null

$wb.Save()
